# Apply the change described by the diff / commit message
# ("Remove unnecessary use of ADL."): delete the
# "Only use ADL where ‘necessary’." bullet, then tidy up the knock-on
# effects of that deletion that Word would otherwise recompute on save:
#   - the trailing lone-space run at the end of the "Investigate whether ..."
#     bullet gets merged back into its preceding text run,
#   - the "_GoBack" bookmark (which tracked the last edit position, right
#     after the "Investigate whether ..." text) moves to the start of the
#     "Optimize all components ..." bullet that now follows the removed
#     bullet directly,
#   - the three "lastRenderedPageBreak" markers (pure layout/repagination
#     bookkeeping) each shift from the bullet that used to start a new page
#     to the next bullet, since the whole list is now one bullet shorter.

$d = $word.ActiveDocument

function Find-ParagraphIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# Replace-ParagraphXml swaps the *entire* contents of the paragraph whose
# text matches $pattern for the supplied raw OOXML - InsertXML replaces
# whatever its target range spans, so the range must cover the complete
# paragraph (all of its runs) or the untouched runs would be lost.
function Replace-ParagraphXml($pattern, $xml) {
    $idx = Find-ParagraphIndex($pattern)
    $p = $d.Paragraphs.Item($idx)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.InsertXML($xml)
}

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. Remove the unnecessary-ADL bullet entirely.
$idx = Find-ParagraphIndex("Only use ADL*")
$d.Paragraphs.Item($idx).Range.Delete()

# 2. Merge the "Investigate whether ..." bullet's separate trailing space
#    run into its main text run.
$investigateIdx = Find-ParagraphIndex("Investigate whether*")
$p = $d.Paragraphs.Item($investigateIdx)
$full = $p.Range.Text
$start = $p.Range.Start
# $full ends with "<space><CR>"; delete just the single space character.
$d.Range($start + $full.Length - 2, $start + $full.Length - 1).Delete()
$p2 = $d.Paragraphs.Item($investigateIdx)
$text2 = $p2.Range.Text
$insertPos = $p2.Range.Start + $text2.Length - 1
$d.Range($insertPos, $insertPos).InsertAfter(" ")

# 3. Relocate the _GoBack bookmark onto the start of the "Optimize all
#    components ..." bullet.
$d.Bookmarks("_GoBack").Delete()
$optimizeIdx = Find-ParagraphIndex("Optimize all components*")
$optStart = $d.Paragraphs.Item($optimizeIdx).Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($optStart, $optStart))

# 4. Shift each lastRenderedPageBreak marker to the next bullet.
Replace-ParagraphXml "New Modules*" "<w:p $w><w:r w:rsidRPr=`"008B06FC`"><w:rPr><w:b/></w:rPr><w:t>New Modules</w:t></w:r></w:p>"
Replace-ParagraphXml "Remote memory*" "<w:p $w><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Remote memory ‘pool’ to avoid allocating entire pages for only a few bytes of data.</w:t></w:r></w:p>"

Replace-ParagraphXml "Transactional hooking*" "<w:p $w><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:t>Transactional hooking.</w:t></w:r></w:p>"
Replace-ParagraphXml "Improved relative instruction rebuilding*" "<w:p $w><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:lastRenderedPageBreak/><w:t>Improved relative instruction rebuilding (including conditionals).</w:t></w:r><w:r w:rsidR=`"00BD30BB`" w:rsidRPr=`"008B06FC`"><w:t xml:space=`"preserve`"> x64 has far more IP relative instructions than x86.</w:t></w:r></w:p>"

Replace-ParagraphXml "Full support for writing back*" "<w:p $w><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:t>Full support for writing back to PE file, including automatically performing adjustments where required to fit in new data or remove unnecessary space.</w:t></w:r></w:p>"
Replace-ParagraphXml "Improve export forwarding*" "<w:p $w><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r w:rsidRPr=`"008B06FC`"><w:lastRenderedPageBreak/><w:t>Improve export forwarding code to detect and handle forward-by-</w:t></w:r><w:r w:rsidR=`"00197559`" w:rsidRPr=`"008B06FC`"><w:t>ordinal</w:t></w:r><w:r w:rsidRPr=`"008B06FC`"><w:t xml:space=`"preserve`"> explicitly rather than forcing the user to detect it and do string manipulation and conversion. </w:t></w:r></w:p>"

Write-Output "Applied: removed unnecessary ADL bullet."
